$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1. All transaction timestamps are now at 0:00 UTC -> row 9's date
#    still carried a fractional day (time-of-day); drop it.
# ------------------------------------------------------------------
$ws1.Range("A9").Value = 44621

# ------------------------------------------------------------------
# 2. Duplicate Sheet1 into a new Sheet2 (placed right after Sheet1)
#    so the extra per-asset breakdown columns can be added without
#    disturbing the original report.
# ------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Name = "Sheet2"

# the copy carried over the trailing blank row (row 10) - drop it
[void]$ws2.Rows("10").Delete()

# ------------------------------------------------------------------
# 3. Sheet2 gets four new header columns (each one repeated further
#    right as a running-total pair), holding the per-asset amounts.
# ------------------------------------------------------------------
$ws2.Range("N1").Value = "ETH Amount"
$ws2.Range("O1").Value = "ETH Value"
$ws2.Range("Q1").Value = "ADA Amount"
$ws2.Range("R1").Value = "ADA Value"
$ws2.Range("T1").Value = "ETH Amount"
$ws2.Range("U1").Value = "ETH Value"
$ws2.Range("W1").Value = "ADA Amount"
$ws2.Range("X1").Value = "ADA Value"

$ws2.Range("N2").Value = 5
$ws2.Range("O2").Value = 5000
$ws2.Range("T2").Value = 5
$ws2.Range("U2").Value = 5000

$ws2.Range("N3").Value = 5
$ws2.Range("O3").Value = 5500
$ws2.Range("T3").Value = 10
$ws2.Range("U3").Value = 10500

$ws2.Range("Q4").Value = 5000
$ws2.Range("R4").Value = 2500
$ws2.Range("W4").Value = 5000
$ws2.Range("X4").Value = 2500

$ws2.Range("Q5").Value = 10000
$ws2.Range("R5").Value = 6000
$ws2.Range("W5").Value = 15000
$ws2.Range("X5").Value = 8500

$ws2.Range("N6").Value = 0.5
$ws2.Range("O6").Value = 1000
$ws2.Range("Q6").Value = 1000
$ws2.Range("T6").Value = 10.5
$ws2.Range("W6").Value = 13000

$ws2.Range("N7").Value = 10
$ws2.Range("Q7").Value = 15200
$ws2.Range("T7").Value = 10
$ws2.Range("W7").Value = 15200

$ws2.Range("N8").Value = 4
$ws2.Range("T8").Value = 4

$ws2.Range("Q9").Value = 2700
$ws2.Range("W9").Value = 2700

# ------------------------------------------------------------------
# 4. Selection bookkeeping to match the edited state: Sheet1's
#    selection becomes the whole data block (no longer the active
#    tab), and Sheet2 - the new active tab - is left with Q6 selected.
# ------------------------------------------------------------------
[void]$ws1.Range("A1:K9").Select()

[void]$ws2.Range("Q6").Select()
